$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the date cell's number format from the row above into the new row's
# date cell, then fill in the new row's values.
$ws.Cells.Item(5, 1).Copy()
$ws.Cells.Item(6, 1).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(6, 1).Value = 42607.889224537037
$ws.Cells.Item(6, 2).Value = 39
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = 0
$ws.Cells.Item(6, 14).Value = "Random"
